$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 keeps its original "quote-prefix text" style (s=1) throughout this script,
# so it's a safe clipboard source whenever we need to re-apply that style to a
# cell after overwriting its value (writing .Value resets a cell's style to
# the plain text style).

# --- Row 6: username 0332743067 -> 02001; password mirrors username; role 2 -> 1 ---
$ws.Range("A6").Value = "02001"
$ws.Range("B6").Value = "02001"
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C6").Value = "1"

# --- Row 7: username 0332743068 -> 02002; password mirrors username; role 2 -> 1 ---
$ws.Range("A7").Value = "02002"
$ws.Range("B7").Value = "02002"
$ws.Range("C7").Value = "1"

# --- Row 8: username 0332743069 -> 02003; password mirrors username; role 2 -> 1 ---
$ws.Range("A8").Value = "02003"
$ws.Range("B8").Value = "02003"
$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C8").Value = "1"

# --- Row 9 (brand-new row): account 02004 ---
$ws.Range("A9").Value = "02004"
$ws.Range("B9").Value = "02004"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "1"

# --- Row 10 (previously just an empty styled placeholder cell): account 02005 ---
$ws.Range("A10").Value = "02005"
$ws.Range("B10").Value = "02005"
$ws.Range("A4").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C10").Value = "1"
$ws.Range("D10").Value = "1"

$excel.CutCopyMode = $false

# Mirror the saved selection/active-cell state.
$ws.Range("D6:D10").Select()
